$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 9 -- shifts old rows 9-17 down to 10-18,
# carrying their formatting (date column D keeps its style) along.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44435
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112052
$ws.Range("G9").Value = "Albahaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 2300
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2400
$ws.Range("N9").Value = "`$/paquete"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 2400
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"

# Append a brand-new row 19 (after the previous last row, now 18) with
# another new weekly record.
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 44432
$ws.Range("D19").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 100112052
$ws.Range("G19").Value = "Albahaca"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 2300
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2400
$ws.Range("N19").Value = "`$/paquete"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 2400
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = "Hortaliza"
